# Insert a new column before column C; the former column C (task descriptions)
# shifts right into column D, leaving a blank column C for a "DONE" status flag.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(3).Insert()

# Mark the first three tasks (rows 3-5) as "DONE" in the new status column C.
$ws.Range("C3").Value() = "DONE"
$ws.Range("C4").Value() = "DONE"
$ws.Range("C5").Value() = "DONE"

# Widen columns B and C to fit the longer content / new status column.
$ws.Range("B:C").ColumnWidth = 27.6

# Move the active selection onto the next task to action (C6).
$ws.Range("C6").Select()
